$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STEPS")

# Insert a new column before column G (7th column), shifting existing
# TC_STEP_ACTION..TC_STEP_CUF_<CODE> columns one to the right.
$ws.Columns("G:G").Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 7).Value = "TC_STEP_CALL_DATASET"

# The newly inserted column keeps the same width as its neighbour (column F)
# but is not re-measured by AutoFit ("bestFit"), just a plain custom width.
$ws.Columns("G:G").ColumnWidth = $ws.Columns("F:F").ColumnWidth

# Keep the active selection on the STEPS sheet at G2, as in the edited file.
$ws.Activate()
$ws.Range("G2").Select()
